# Weekly CompStat data refresh: new crime data collected for the period
# ending 12/10/2023 (Volume 30, Number 49).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header: volume/issue number and the reporting week date range ---
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Crime-statistics grid (rows 14-29): updated counts and recomputed % changes ---
$ws.Range("M14").Value = -21.428571428571
$ws.Range("N14").Value = -82.8125
$ws.Range("G15").Value = 2
$ws.Range("I15").Value = 39
$ws.Range("K15").Value = -13.333333333333
$ws.Range("L15").Value = -7.142857142857
$ws.Range("M15").Value = 18.181818181818
$ws.Range("N15").Value = -42.647058823529
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = 44
$ws.Range("H16").Value = -9.090909090909
$ws.Range("I16").Value = 507
$ws.Range("J16").Value = 698
$ws.Range("K16").Value = -27.363896848137
$ws.Range("L16").Value = 3.680981595092
$ws.Range("M16").Value = 3.048780487804
$ws.Range("N16").Value = -73.959938366718
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = -52.173913043478
$ws.Range("F17").Value = 56
$ws.Range("G17").Value = 70
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 795
$ws.Range("J17").Value = 765
$ws.Range("K17").Value = 3.92156862745
$ws.Range("L17").Value = 18.834080717488
$ws.Range("M17").Value = 61.914460285132
$ws.Range("N17").Value = -11.371237458194
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 85.714285714285
$ws.Range("I18").Value = 308
$ws.Range("J18").Value = 288
$ws.Range("K18").Value = 6.944444444444
$ws.Range("L18").Value = 32.758620689655
$ws.Range("M18").Value = -16.981132075471
$ws.Range("N18").Value = -83.178590933915
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 27
$ws.Range("E19").Value = -22.222222222222
$ws.Range("G19").Value = 94
$ws.Range("H19").Value = -22.340425531914
$ws.Range("I19").Value = 874
$ws.Range("J19").Value = 994
$ws.Range("K19").Value = -12.072434607645
$ws.Range("L19").Value = 15.455746367239
$ws.Range("M19").Value = 61.552680221811
$ws.Range("N19").Value = 17.789757412398
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 21
$ws.Range("E20").Value = -52.380952380952
$ws.Range("F20").Value = 49
$ws.Range("G20").Value = 62
$ws.Range("H20").Value = -20.967741935483
$ws.Range("I20").Value = 615
$ws.Range("J20").Value = 520
$ws.Range("K20").Value = 18.26923076923
$ws.Range("L20").Value = 46.778042959427
$ws.Range("M20").Value = 168.558951965065
$ws.Range("N20").Value = -68.266253869969
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 89
$ws.Range("E21").Value = -37.078651685393
$ws.Range("F21").Value = 244
$ws.Range("G21").Value = 286
$ws.Range("H21").Value = -14.685314685314
$ws.Range("I21").Value = 3149
$ws.Range("J21").Value = 3316
$ws.Range("K21").Value = -5.036188178528
$ws.Range("L21").Value = 20.007621951219
$ws.Range("M21").Value = 45.048364808843
$ws.Range("N21").Value = -57.940430078803
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 15
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -11.764705882352
$ws.Range("M22").Value = 7.142857142857
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 14
$ws.Range("H23").Value = -54.838709677419
$ws.Range("I23").Value = 287
$ws.Range("J23").Value = 320
$ws.Range("K23").Value = -10.3125
$ws.Range("L23").Value = 13.438735177865
$ws.Range("M23").Value = 32.258064516129
$ws.Range("C24").Value = 41
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = 7.894736842105
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 168
$ws.Range("H24").Value = -29.166666666666
$ws.Range("I24").Value = 1814
$ws.Range("J24").Value = 1967
$ws.Range("K24").Value = -7.778342653787
$ws.Range("L24").Value = 31.259044862518
$ws.Range("M24").Value = 25.189786059351
$ws.Range("C25").Value = 32
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 77.777777777777
$ws.Range("F25").Value = 104
$ws.Range("G25").Value = 74
$ws.Range("H25").Value = 40.54054054054
$ws.Range("I25").Value = 1069
$ws.Range("J25").Value = 1062
$ws.Range("K25").Value = 0.659133709981
$ws.Range("L25").Value = 12.882787750792
$ws.Range("M25").Value = -27.229407760381
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 58
$ws.Range("K26").Value = -12.121212121212
$ws.Range("L26").Value = -20.547945205479
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 69
$ws.Range("K27").Value = 43.478260869565
$ws.Range("L27").Value = 1.020408163265
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "***.*"
$ws.Range("M28").Value = -21.56862745098
$ws.Range("N28").Value = -73.154362416107
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "***.*"
$ws.Range("M29").Value = -20.930232558139
$ws.Range("N29").Value = -74.814814814814

# --- A few cells switched between numeric and text ("0" / "***.*") representation;
#     copy number-format/style from an unaffected sibling cell so the style matches. ---
$ws.Range("J22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("I27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("I27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("L27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Copy() | Out-Null
$ws.Range("G28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Copy() | Out-Null
$ws.Range("H28").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

